$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header rename
$ws.Range("O1").Value = "F1 train"

# Row 2
$ws.Range("O2").Value = 0.8974358974358975

# Row 3
$ws.Range("O3").Value = 1

# Row 4
$ws.Range("O4").Value = 1

# Row 5
$ws.Range("O5").Value = 0.8205128205128205

# Row 6 (MLP, hidden_layer_sizes (32,) -> (64,), activation relu -> tanh)
$ws.Range("C6").Value = "{'activation': 'tanh', 'alpha': 0.0001, 'hidden_layer_sizes': (64,), 'learning_rate': 'constant'}"
$ws.Range("E6").Value = 5
$ws.Range("F6").Value = 4
$ws.Range("G6").Value = 6
$ws.Range("H6").Value = 5
$ws.Range("I6").Value = 0.55
$ws.Range("J6").Value = 0.5263157894736842
$ws.Range("K6").Value = 0.5
$ws.Range("L6").Value = 0.5555555555555556
$ws.Range("M6").Value = 0.6
$ws.Range("N6").Value = 0.5
$ws.Range("O6").Value = 0.7560975609756098

# Row 7
$ws.Range("O7").Value = 0.9866666666666667

# Row 8
$ws.Range("O8").Value = 1

# Row 9
$ws.Range("O9").Value = 0.9736842105263158

# Row 10
$ws.Range("O10").Value = 0.7073170731707317

# Row 11
$ws.Range("O11").Value = 0.6440677966101694

# Row 12
$ws.Range("O12").Value = 0.9210526315789473

# Row 13
$ws.Range("O13").Value = 0.96

# Row 14
$ws.Range("O14").Value = 0.821917808219178

# Row 15
$ws.Range("O15").Value = 0.7142857142857143

# Row 16 (MLP, hidden_layer_sizes (32, 16) -> (64, 32), activation relu -> tanh)
$ws.Range("C16").Value = "{'activation': 'tanh', 'alpha': 0.0001, 'hidden_layer_sizes': (64, 32), 'learning_rate': 'constant'}"
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 5
$ws.Range("G16").Value = 5
$ws.Range("H16").Value = 7
$ws.Range("I16").Value = 0.4
$ws.Range("J16").Value = 0.3333333333333333
$ws.Range("K16").Value = 0.3
$ws.Range("L16").Value = 0.375
$ws.Range("M16").Value = 0.5
$ws.Range("N16").Value = 0.3
$ws.Range("O16").Value = 0.7435897435897436
